# "Working with Multiple Sheets"
#
# The name/colour columns that used to live on Sheet1 move over to
# Sheet2, and Sheet1 keeps just name / price_used / quantity (the
# colour column is removed so price_used and quantity shift left).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Copy the name (A) and colour (B) columns, header included, over to Sheet2.
[void]$ws1.Range("A1:B9").Copy($ws2.Range("A1"))

# Remove the colour column from Sheet1; price_used/quantity shift left
# into columns B/C.
[void]$ws1.Range("B1:B9").EntireColumn.Delete()

# Match the final selections recorded in each sheet view: select the
# copied block on Sheet2 first, then return focus to Sheet1 so it stays
# the active tab with its own new selection.
[void]$ws2.Range("A1:B9").Select()
[void]$ws1.Range("F8").Select()
